# Update totaalstand (total standings) sheet with new match results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-13: Rang(A), Speler(B), Score(C), 180'ers(D), 100+ finishes(E),
# Totaal Score(F), Aantal Darts(G), 3-Darts Gemiddelde(H), Totaal(I), Winnaar(J)
$data = @(
    @(1,  "Burger Peach",            29, 0, 2, 17293, 796, 65.17,              31, 2),
    @(2,  "Yannick den Daggelder",   20, 0, 0, 11319, 619, 54.86,              20, 1),
    @(3,  "Niels van Dommelen",      17, 0, 0, 13810, 794, 52.18,              17, 0),
    @(4,  "Nick Fitzpatrick",        14, 1, 0, 8568,  401, 64.09999999999999,  15, 0),
    @(5,  "Rocky Van Den Eeckhoudt", 13, 0, 1, 11777, 619, 57.08,              14, 0),
    @(6,  "Sion Foulkes",             5, 0, 0, 5404,  367, 44.17,               5, 0),
    @(6,  "Lukas G",                  4, 1, 0, 5294,  283, 56.12,               5, 0),
    @(8,  "Constantinos Mavroudis",   4, 0, 0, 1991,  134, 44.57,               4, 0),
    @(8,  "Nigel Riedel",             3, 1, 0, 4425,  249, 53.31,               4, 0),
    @(10, "Danny Littler",            3, 0, 0, 0,     0,   $null,               3, 0),
    @(10, "Noah B",                   3, 0, 0, 4473,  301, 44.58,               3, 0),
    @(12, "Quintin Marais",           2, 0, 0, 2388,  107, 66.95,               2, 0)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    if ($null -eq $rec[7]) {
        $ws.Cells.Item($row, 8).Value = ""
    } else {
        $ws.Cells.Item($row, 8).Value = $rec[7]
    }
    $ws.Cells.Item($row, 9).Value = $rec[8]
    $ws.Cells.Item($row, 10).Value = $rec[9]
    $row = $row + 1
}
